$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, copying formatting from G1 (the neighboring header cell)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$saveValues = @(0, 0, 0, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
